$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with refreshed crypto data.
# D-column values are forced to text (NumberFormat "@") so that Excel does not
# auto-convert numeric-looking strings (losing trailing zeros / turning them into
# real numbers); the style is then reset to "Normal" so no new cell style lingers.

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '27.365.36'
$c.Style = "Normal"

$ws.Range("E2").Value = '  +9.32%  '

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '1.788.32'
$c.Style = "Normal"

$ws.Range("E3").Value = '  +6.66%  '

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '1.002'
$c.Style = "Normal"

$ws.Range("E4").Value = '  +0.06%  '

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '339.25'
$c.Style = "Normal"

$ws.Range("E5").Value = '  +2.84%  '

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '0.9974'
$c.Style = "Normal"

$ws.Range("E6").Value = '  -0.14%  '

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.3800'
$c.Style = "Normal"

$ws.Range("E7").Value = '  +3.99%  '

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.3496'
$c.Style = "Normal"

$ws.Range("E8").Value = '  +7.77%  '

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '49.58'
$c.Style = "Normal"

$ws.Range("E9").Value = '  +4.72%  '

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '1.220'
$c.Style = "Normal"

$ws.Range("E10").Value = '  +5.92%  '

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.07641'
$c.Style = "Normal"

$ws.Range("E11").Value = '  +6.04%  '

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.9973'
$c.Style = "Normal"

$ws.Range("E12").Value = '  -0.23%  '

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '6.597'
$c.Style = "Normal"

$ws.Range("E13").Value = '  +8.17%  '

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '21.54'
$c.Style = "Normal"

$ws.Range("E14").Value = '  +9.18%  '

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '7.211'
$c.Style = "Normal"

$ws.Range("E15").Value = '  +8.09%  '

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '1.786.28'
$c.Style = "Normal"

$ws.Range("E16").Value = '  +6.66%  '

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '0.00001119'
$c.Style = "Normal"

$ws.Range("E17").Value = '  +6.05%  '

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '0.06773'
$c.Style = "Normal"

$ws.Range("E18").Value = '  +3.76%  '

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '85.22'
$c.Style = "Normal"

$ws.Range("E19").Value = '  +7.81%  '

$ws.Range("E20").Value = '  -0.09%  '

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '17.67'
$c.Style = "Normal"

$ws.Range("E21").Value = '  +11.33%  '

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '6.405'
$c.Style = "Normal"

$ws.Range("E22").Value = '  +8.09%  '

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '13.19'
$c.Style = "Normal"

$ws.Range("E23").Value = '  +2.49%  '

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '27.380.18'
$c.Style = "Normal"

$ws.Range("E24").Value = '  +9.39%  '

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '2.473'
$c.Style = "Normal"

$ws.Range("E25").Value = '  +1.46%  '

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '1.541'
$c.Style = "Normal"

$ws.Range("E26").Value = '  +28.75%  '

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '2.550'
$c.Style = "Normal"

$ws.Range("E27").Value = '  +6.68%  '

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '20.28'
$c.Style = "Normal"

$ws.Range("E28").Value = '  +8.13%  '

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '153.47'
$c.Style = "Normal"

$ws.Range("E29").Value = '  +2.79%  '

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '1.984.64'
$c.Style = "Normal"

$ws.Range("E30").Value = '  +6.70%  '

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '135.69'
$c.Style = "Normal"

$ws.Range("E31").Value = '  +7.59%  '

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '6.408'
$c.Style = "Normal"

$ws.Range("E32").Value = '  +9.93%  '

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '4.125'
$c.Style = "Normal"

$ws.Range("E33").Value = '  +0.71%  '

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '0.08754'
$c.Style = "Normal"

$ws.Range("E34").Value = '  +3.32%  '

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '13.45'
$c.Style = "Normal"

$ws.Range("E35").Value = '  +8.51%  '

$ws.Range("E36").Value = '  +3.00%  '

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '5.639'
$c.Style = "Normal"

$ws.Range("E37").Value = '  +8.88%  '

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.02421'
$c.Style = "Normal"

$ws.Range("E38").Value = '  +8.01%  '

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.2265'
$c.Style = "Normal"

$ws.Range("E39").Value = '  +8.20%  '

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '0.06501'
$c.Style = "Normal"

$ws.Range("E40").Value = '  +6.49%  '

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '0.6776'
$c.Style = "Normal"

$ws.Range("E41").Value = '  +13.15%  '

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '8.802'
$c.Style = "Normal"

$ws.Range("E42").Value = '  +6.21%  '

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '1.242'
$c.Style = "Normal"

$ws.Range("E43").Value = '  +0.52%  '

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '14.80'
$c.Style = "Normal"

$ws.Range("E44").Value = '  +8.61%  '

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '0.6463'
$c.Style = "Normal"

$ws.Range("E45").Value = '  +12.67%  '

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '0.9977'
$c.Style = "Normal"

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '3.979'
$c.Style = "Normal"

$ws.Range("E47").Value = '  +3.83%  '

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '2.158'
$c.Style = "Normal"

$ws.Range("E48").Value = '  +9.56%  '

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '131.68'
$c.Style = "Normal"

$ws.Range("E49").Value = '  +5.76%  '

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '0.07361'
$c.Style = "Normal"

$ws.Range("E50").Value = '  +4.85%  '

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '80.34'
$c.Style = "Normal"

$ws.Range("E51").Value = '  +7.17%  '
